$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 208, pushing existing rows 208-312
# down to 210-314 (matches the shift seen across the whole diff).
$ws.Rows.Item(208).Insert()
$ws.Rows.Item(209).Insert()

# New row 208: Ajo / Chino / Primera, 2022-09-02 (date serial 44806)
$ws.Cells.Item(208, 1).Value = 4
$ws.Cells.Item(208, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(208, 3).Value = "Los Lagos"
$ws.Cells.Item(208, 4).Value = 44806
$ws.Cells.Item(208, 5).Value = 10
$ws.Cells.Item(208, 6).Value = 100112003
$ws.Cells.Item(208, 7).Value = "Ajo"
$ws.Cells.Item(208, 8).Value = "Chino"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 120
$ws.Cells.Item(208, 11).Value = 28000
$ws.Cells.Item(208, 12).Value = 28000
$ws.Cells.Item(208, 13).Value = 28000
$ws.Cells.Item(208, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(208, 15).Value = "China"
$ws.Cells.Item(208, 16).Value = 2800
$ws.Cells.Item(208, 17).Value = 10
$ws.Cells.Item(208, 18).Value = "Hortaliza"

# New row 209: Ajo / Chino / Segunda, same date (serial 44806)
$ws.Cells.Item(209, 1).Value = 4
$ws.Cells.Item(209, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(209, 3).Value = "Los Lagos"
$ws.Cells.Item(209, 4).Value = 44806
$ws.Cells.Item(209, 5).Value = 10
$ws.Cells.Item(209, 6).Value = 100112003
$ws.Cells.Item(209, 7).Value = "Ajo"
$ws.Cells.Item(209, 8).Value = "Chino"
$ws.Cells.Item(209, 9).Value = "Segunda"
$ws.Cells.Item(209, 10).Value = 120
$ws.Cells.Item(209, 11).Value = 26000
$ws.Cells.Item(209, 12).Value = 26000
$ws.Cells.Item(209, 13).Value = 26000
$ws.Cells.Item(209, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(209, 15).Value = "China"
$ws.Cells.Item(209, 16).Value = 2600
$ws.Cells.Item(209, 17).Value = 10
$ws.Cells.Item(209, 18).Value = "Hortaliza"
